$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so values like "1.00" or "571.69"
# are written as literal text (matching the source inlineStr cells) instead of
# being auto-converted to numbers by Excel's type inference.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "65.056.27"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "3.146.51"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "571.69"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").Value = "150.07"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.143.22"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +6.16%  "
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "0.506"
$ws.Range("E12").Value = "  +7.67%  "
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  +11.79%  "
$ws.Range("E14").Value = "  +6.52%  "
$ws.Range("D15").Value = "3.662.53"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "65.080.29"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "7.22"
$ws.Range("E17").Value = "  +6.58%  "
$ws.Range("D18").Value = "3.147.55"
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "510.62"
$ws.Range("E20").Value = "  +6.71%  "
$ws.Range("D21").Value = "14.91"
$ws.Range("E21").Value = "  +7.08%  "
$ws.Range("D22").Value = "0.731"
$ws.Range("E22").Value = "  +8.24%  "
$ws.Range("D23").Value = "15.57"
$ws.Range("E23").Value = "  +13.99%  "
$ws.Range("D24").Value = "7.85"
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("D25").Value = "85.61"
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E27").Value = "  +4.59%  "
$ws.Range("D28").Value = "8.76"
$ws.Range("E28").Value = "  +8.05%  "
$ws.Range("E29").Value = "  +5.24%  "
$ws.Range("D30").Value = "28.02"
$ws.Range("E30").Value = "  +7.03%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("E33").Value = "  +6.15%  "
$ws.Range("D34").Value = "6.05"
$ws.Range("E34").Value = "  +8.62%  "
$ws.Range("D35").Value = "6.62"
$ws.Range("E35").Value = "  +7.02%  "
$ws.Range("D36").Value = "55.64"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Value = "472.93"
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("E38").Value = "  +4.50%  "
$ws.Range("D39").Value = "0.0858"
$ws.Range("E39").Value = "  +3.65%  "
$ws.Range("D40").Value = "3.03"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").Value = "3.126.81"
$ws.Range("E41").Value = "  +5.28%  "
$ws.Range("E42").Value = "  +4.66%  "
$ws.Range("D43").Value = "0.119"
$ws.Range("E43").Value = "  +3.94%  "
$ws.Range("E44").Value = "  +11.02%  "
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  +12.92%  "
$ws.Range("D46").Value = "29.30"
$ws.Range("E46").Value = "  +5.21%  "
$ws.Range("D47").Value = "0.0₃0575"
$ws.Range("E47").Value = "  +11.57%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  +4.01%  "
$ws.Range("E50").Value = "  +11.62%  "
$ws.Range("D51").Value = "118.41"
$ws.Range("E51").Value = "  -1.85%  "

# Restore default (General) formatting so no stray number-format styles remain
$rng.ClearFormats()

